# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The whole-table refresh timestamp moved from serial 45171 (2023-09-02)
# to serial 45172 (2023-09-03) for every record, rows 2 through 265.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C265").Value = 45172
